# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" / "Error Detail" columns for the 0f0b2084... row (row 5) on both the
# zh-cn and de-de sheets, because the handback file that came back was not built
# against the latest handoff -- so we now link back to the target markdown file,
# record the (new) handback xlf + timestamp, and surface the version-mismatch
# error message. Also widens the J/K/R columns so the new content is readable.

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/be6ba416f379ac23bbbae3a163d2faf055f8b80c/e2e/0f0b2084-9fd7-447d-bc04-b237989b967e.md"
$targetDisplay = "0f0b2084-9fd7-447d-bc04-b237989b967e.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/56e4c8eb81fe7d53ce3ec971a445a151ed221231/e2e/0f0b2084-9fd7-447d-bc04-b237989b967e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/be6ba416f379ac23bbbae3a163d2faf055f8b80c/e2e/0f0b2084-9fd7-447d-bc04-b237989b967e.md."

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Columns.Item(10).ColumnWidth = 39.1666667
$wsZh.Columns.Item(11).ColumnWidth = 39.1666667
$wsZh.Columns.Item(18).ColumnWidth = 39.1666667

$wsZh.Range("K5").Value = "0f0b2084-9fd7-447d-bc04-b237989b967e.fcf151572a7902e77872e64eff9d21ac69f842e5.zh-cn.xlf"
$wsZh.Range("L5").Value = "2017-02-21 09:26:23"
$wsZh.Range("R5").Value = $errorDetail
$wsZh.Hyperlinks.Add($wsZh.Range("J5"), $targetUrl, "", "", $targetDisplay)

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(10).ColumnWidth = 39.1666667
$wsDe.Columns.Item(11).ColumnWidth = 39.1666667
$wsDe.Columns.Item(18).ColumnWidth = 39.1666667

$wsDe.Range("K5").Value = "0f0b2084-9fd7-447d-bc04-b237989b967e.fcf151572a7902e77872e64eff9d21ac69f842e5.de-de.xlf"
$wsDe.Range("L5").Value = "2017-02-21 09:26:47"
$wsDe.Range("R5").Value = $errorDetail
$wsDe.Hyperlinks.Add($wsDe.Range("J5"), $targetUrl, "", "", $targetDisplay)
